$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old cells that are no longer used
$ws.Range("B5").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B13").ClearContents()

# Set the new values in the order that matches the target shared-string table
$ws.Range("B13").Value = "Wahl des Nachbardreiecks in Promenade"
$ws.Range("B11").Value = "Visualisierung -> Pfad in Maillage"
$ws.Range("B9").Value = "Listen Ansatz (sort)"
$ws.Range("B15").Value = "Laufzeit / Komplexität"
$ws.Range("B17").Value = "Promenade falls kein Nachbar existiert (-1)"

# Update selection to match the target state
$ws.Range("B17").Select()
